$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'58.017.72"
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +3.02%  '
$ws.Range('D3').Value = "'3.075.55"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +2.30%  '
$ws.Range('D5').Value = "'518.25"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.82%  '
$ws.Range('D6').Value = "'143.07"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +3.18%  '
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('E8').Value = '  +1.67%  '
$ws.Range('D9').Value = "'7.31"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +2.55%  '
$ws.Range('E10').Value = '  +0.44%  '
$ws.Range('E11').Value = '  +2.65%  '
$ws.Range('D12').Value = "'3.599.82"
$ws.Range('D12').ClearFormats()
$ws.Range('E13').Value = '  +3.03%  '
$ws.Range('E14').Value = '  +3.41%  '
$ws.Range('E15').Value = '  +1.66%  '
$ws.Range('D16').Value = "'58.027.99"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = "'3.076.16"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +2.24%  '
$ws.Range('D18').Value = "'6.11"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.64%  '
$ws.Range('D19').Value = "'12.92"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('E20').Value = '  +1.48%  '
$ws.Range('D21').Value = "'333.72"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.17%  '
$ws.Range('E22').Value = '  +0.01%  '
$ws.Range('D23').Value = "'0.502"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('D24').Value = "'65.69"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('E25').Value = '  +2.77%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  -2.86%  '
$ws.Range('D28').Value = "'6.50"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +2.57%  '
$ws.Range('D29').Value = "'7.28"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +6.00%  '
$ws.Range('E30').Value = '  +1.95%  '
$ws.Range('D31').Value = "'1.20"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +3.29%  '
$ws.Range('D32').Value = "'20.78"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.97%  '
$ws.Range('D33').Value = "'154.70"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.24%  '
$ws.Range('D34').Value = "'4.58"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.19%  '
$ws.Range('D35').Value = "'6.02"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.40%  '
$ws.Range('D36').Value = "'27.04"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.61%  '
$ws.Range('E37').Value = '  +4.33%  '
$ws.Range('D38').Value = "'0.0678"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.57%  '
$ws.Range('D39').Value = "'3.114.14"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +2.30%  '
$ws.Range('E40').Value = '  +3.87%  '
$ws.Range('D41').Value = "'36.48"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.31%  '
$ws.Range('D42').Value = "'1.00"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +0.01%  '
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = "'2.269.51"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.12%  '
$ws.Range('E45').Value = '  +8.22%  '
$ws.Range('D46').Value = "'21.10"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +8.36%  '
$ws.Range('D47').Value = "'1.38"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.66%  '
$ws.Range('D48').Value = "'0.948"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +2.82%  '
$ws.Range('D49').Value = "'5.94"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +1.69%  '
$ws.Range('D50').Value = "'0.748"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +11.26%  '
$ws.Range('D51').Value = "'257.23"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +12.93%  '
